$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ D = "282.52"; E = "1.64%"; G = "10" }
    3 = @{ D = "28.38"; E = "4.03%"; G = "10" }
    4 = @{ D = "5.047"; E = "3.52%"; G = "10" }
    5 = @{ D = "0.06497"; E = "0.97%"; G = "10" }
    6 = @{ D = "7.214"; E = "2.91%"; G = "10" }
    7 = @{ D = "1.393"; E = "17.32%"; G = "10" }
    8 = @{ D = "0.9181"; E = "3.57%"; G = "10" }
    9 = @{ D = "0.1532"; E = "-1.82%"; G = "10" }
    10 = @{ D = "0.06383"; E = "25.01%"; G = "10" }
    11 = @{ D = "0.07620"; E = "1.67%"; G = "10" }
    12 = @{ D = "0.02824"; E = "-2.27%"; G = "10" }
    13 = @{ D = "0.08966"; E = "-0.11%"; G = "10" }
    14 = @{ D = "0.001588"; E = "1.19%"; G = "10" }
    15 = @{ D = "0.0006344"; E = "-0.72%"; G = "10" }
    16 = @{ D = "0.006112"; E = "-0.64%"; G = "10" }
    17 = @{ E = "-1.02%"; G = "10" }
    18 = @{ E = "1.60%"; G = "10" }
    19 = @{ E = "-1.40%"; G = "10" }
    20 = @{ G = "10" }
    21 = @{ E = "-0.71%"; G = "10" }
    22 = @{ D = "3.976"; E = "0.94%"; G = "10" }
    23 = @{ G = "10" }
    24 = @{ D = "0.04438"; E = "0.55%"; G = "10" }
    25 = @{ D = "0.001183"; E = "0.56%"; G = "10" }
    26 = @{ D = "0.004457"; E = "15.04%"; G = "10" }
    27 = @{ G = "10" }
    28 = @{ E = "1.71%"; G = "10" }
    29 = @{ E = "-1.97%"; G = "10" }
    30 = @{ G = "10" }
    31 = @{ G = "10" }
    32 = @{ G = "10" }
    33 = @{ G = "10" }
    34 = @{ G = "10" }
    35 = @{ G = "10" }
    36 = @{ G = "10" }
    37 = @{ G = "10" }
    38 = @{ G = "10" }
    39 = @{ G = "10" }
    40 = @{ D = "0.04099"; E = "-0.99%"; G = "10" }
    41 = @{ D = "0.006682"; E = "-1.14%"; G = "10" }
    42 = @{ D = "0.1232"; E = "4.95%"; G = "10" }
    43 = @{ D = "0.002150"; E = "14.38%"; G = "10" }
    44 = @{ D = "0.01150"; E = "2.49%"; G = "10" }
    45 = @{ D = "0.00005395"; E = "1.51%"; G = "10" }
    46 = @{ D = "2.000"; E = "18.53%"; G = "10" }
    47 = @{ E = "-0.12%"; G = "10" }
    48 = @{ G = "10" }
    49 = @{ G = "10" }
    50 = @{ G = "10" }
    51 = @{ G = "10" }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    foreach ($col in $vals.Keys) {
        $addr = "$col$row"
        $cell = $ws.Range($addr)
        $cell.NumberFormat = "@"
        $cell.Value = $vals[$col]
        $cell.Style = "Normal"
    }
}